# Auto-generated edit script applying numeric updates to the Sheets workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 7303
$ws.Range("I31").Value = 7303
$ws.Range("K31").Value = 21909
$ws.Range("M31").Value = -21679
$ws.Range("H40").Value = 5023.205
$ws.Range("I40").Value = 2960
$ws.Range("J40").Value = 5259
$ws.Range("K40").Value = 2960
$ws.Range("L40").Value = 5259
$ws.Range("M40").Value = -2785
$ws.Range("N40").Value = -5609
$ws.Range("H62").Value = 3035.524
$ws.Range("I62").Value = 3035.524
$ws.Range("K62").Value = 3035.524
$ws.Range("M62").Value = -2411.524
$ws.Range("H65").Value = 3035.524
$ws.Range("I65").Value = 3035.524
$ws.Range("K65").Value = 15177.62
$ws.Range("M65").Value = -12057.62
$ws.Range("H100").Value = 3849.6897
$ws.Range("I100").Value = 1424.4546
$ws.Range("J100").Value = 5331.778
$ws.Range("K100").Value = 1424.4546
$ws.Range("L100").Value = 5331.778
$ws.Range("M100").Value = -883.4546
$ws.Range("N100").Value = -6413.778
$ws.Range("H103").Value = 720.44446
$ws.Range("I103").Value = 502.66666
$ws.Range("J103").Value = 829.3333
$ws.Range("K103").Value = 1507.99998
$ws.Range("L103").Value = 2487.9999
$ws.Range("M103").Value = -921.9999800000001
$ws.Range("N103").Value = -3659.9999
$ws.Range("H135").Value = 4529.95
$ws.Range("J135").Value = 6995.375
$ws.Range("L135").Value = 62958.375
$ws.Range("N135").Value = -68028.375
$ws.Range("H138").Value = 3528.8276
$ws.Range("I138").Value = 2340.2354
$ws.Range("K138").Value = 7020.706200000001
$ws.Range("M138").Value = -1880.706200000001
$ws.Range("H141").Value = 6194.048
$ws.Range("I141").Value = 2906.8333
$ws.Range("J141").Value = 10577
$ws.Range("K141").Value = 8720.499899999999
$ws.Range("L141").Value = 31731
$ws.Range("M141").Value = -3540.499899999999
$ws.Range("N141").Value = -42091
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2333.1853
$ws.Range("I61").Value = 2309.1365
$ws.Range("K61").Value = 2309.1365
$ws.Range("M61").Value = -2097.1365
$ws.Range("H74").Value = 5732.8867
$ws.Range("J74").Value = 18151.3
$ws.Range("L74").Value = 18151.3
$ws.Range("N74").Value = -19899.3
$ws.Range("H77").Value = 5732.8867
$ws.Range("J77").Value = 18151.3
$ws.Range("L77").Value = 90756.5
$ws.Range("N77").Value = -99492.5
$ws.Range("H110").Value = 1089.7727
$ws.Range("I110").Value = 935.6429000000001
$ws.Range("K110").Value = 935.6429000000001
$ws.Range("M110").Value = 1109.3571
$ws.Range("H132").Value = 3824.4138
$ws.Range("J132").Value = 6822.448
$ws.Range("L132").Value = 20467.344
$ws.Range("N132").Value = -25527.344
$ws.Range("H136").Value = 2333.1853
$ws.Range("I136").Value = 2309.1365
$ws.Range("K136").Value = 6927.4095
$ws.Range("M136").Value = -4377.4095
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 71999.5
$ws.Range("J51").Value = 71999.5
$ws.Range("L51").Value = 71999.5
$ws.Range("N51").Value = -72981.5
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H82").Value = 31106.143
$ws.Range("J82").Value = 48333.332
$ws.Range("L82").Value = 48333.332
$ws.Range("N82").Value = -49099.332
$ws.Range("H85").Value = 31106.143
$ws.Range("J85").Value = 48333.332
$ws.Range("L85").Value = 48333.332
$ws.Range("N85").Value = -50985.332
$ws.Range("H94").Value = 4082.0715
$ws.Range("I94").Value = 4339.273
$ws.Range("K94").Value = 4339.273
$ws.Range("M94").Value = -3888.273
$ws.Range("H134").Value = 2380.8918
$ws.Range("I134").Value = 2131.8845
$ws.Range("K134").Value = 6395.6535
$ws.Range("M134").Value = -3860.6535
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5922.3667
$ws.Range("J58").Value = 7971.6665
$ws.Range("L58").Value = 7971.6665
$ws.Range("N58").Value = -8377.666499999999
$ws.Range("H68").Value = 42200
$ws.Range("J68").Value = 42200
$ws.Range("L68").Value = 42200
$ws.Range("N68").Value = -43698
$ws.Range("H71").Value = 42200
$ws.Range("J71").Value = 42200
$ws.Range("L71").Value = 126600
$ws.Range("N71").Value = -134088
$ws.Range("H107").Value = 1223.16
$ws.Range("I107").Value = 932.17645
$ws.Range("J107").Value = 1841.5
$ws.Range("K107").Value = 932.17645
$ws.Range("L107").Value = 1841.5
$ws.Range("M107").Value = 987.82355
$ws.Range("N107").Value = -5681.5
$ws.Range("H122").Value = 10379.226
$ws.Range("J122").Value = 44566.668
$ws.Range("L122").Value = 133700.004
$ws.Range("N122").Value = -138600.004
$ws.Range("H136").Value = 5922.3667
$ws.Range("J136").Value = 7971.6665
$ws.Range("L136").Value = 23914.9995
$ws.Range("N136").Value = -29014.9995
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1151.3529
$ws.Range("I5").Value = 985.55554
$ws.Range("K5").Value = 2956.66662
$ws.Range("M5").Value = -2844.66662
$ws.Range("H86").Value = 239.25
$ws.Range("I86").Value = 315.6
$ws.Range("K86").Value = 946.8000000000001
$ws.Range("M86").Value = 239.1999999999999
$ws.Range("H89").Value = 239.25
$ws.Range("I89").Value = 315.6
$ws.Range("K89").Value = 2840.4
$ws.Range("M89").Value = 3087.6
$ws.Range("H104").Value = 5137
$ws.Range("J104").Value = 4895
$ws.Range("L104").Value = 14685
$ws.Range("N104").Value = -19927
$ws.Range("H107").Value = 3324.1875
$ws.Range("I107").Value = 2290
$ws.Range("J107").Value = 3562.8462
$ws.Range("K107").Value = 6870
$ws.Range("L107").Value = 10688.5386
$ws.Range("M107").Value = -4950
$ws.Range("N107").Value = -14528.5386
$ws.Range("H122").Value = 1614473.5
$ws.Range("J122").Value = 1906.75
$ws.Range("L122").Value = 17160.75
$ws.Range("N122").Value = -22060.75
$ws.Range("H135").Value = 1151.3529
$ws.Range("I135").Value = 985.55554
$ws.Range("K135").Value = 8869.99986
$ws.Range("M135").Value = -6334.99986
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3024.8
$ws.Range("I102").Value = 3096.7693
$ws.Range("K102").Value = 3096.7693
$ws.Range("M102").Value = -1474.7693
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 1995
$ws.Range("I57").Value = 1995
$ws.Range("K57").Value = 1995
$ws.Range("M57").Value = -1429
$ws.Range("H104").Value = 20000
$ws.Range("J104").Value = 20000
$ws.Range("L104").Value = 20000
$ws.Range("N104").Value = -26988
$ws.Range("H131").Value = 53800
$ws.Range("J131").Value = 53800
$ws.Range("L131").Value = 53800
$ws.Range("N131").Value = -63880
$ws.Range("H136").Value = 11005.25
$ws.Range("I136").Value = 4910
$ws.Range("J136").Value = 14052.875
$ws.Range("K136").Value = 14730
$ws.Range("L136").Value = 42158.625
$ws.Range("M136").Value = -12180
$ws.Range("N136").Value = -47258.625
$ws.Range("H138").Value = 10000
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 128414.22
$ws.Range("I122").Value = 3358.1428
$ws.Range("K122").Value = 10074.4284
$ws.Range("M122").Value = -7624.428400000001
$ws.Range("H132").Value = 2735.9023
$ws.Range("J132").Value = 3679.5454
$ws.Range("L132").Value = 11038.6362
$ws.Range("N132").Value = -16098.6362
$ws.Range("H136").Value = 2182.7856
$ws.Range("J136").Value = 2051.5
$ws.Range("L136").Value = 6154.5
$ws.Range("N136").Value = -11254.5
